$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.557.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "'1.885.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("D5").Value = "'244.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.17%  "
$ws.Range("E6").Value = "  +2.51%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "'42.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.29%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "'0.0706"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'2.155.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "'12.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.01%  "
$ws.Range("D14").Value = "'1.928.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.46%  "
$ws.Range("D15").Value = "'0.689"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "'4.80"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").Value = "'35.516.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "'71.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").Value = "'0.0₃0810"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "'243.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "'4.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "'170.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").Value = "'2.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.27%  "
$ws.Range("D27").Value = "'8.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.31%  "
$ws.Range("D28").Value = "'17.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").Value = "'0.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +30.79%  "
$ws.Range("D31").Value = "'0.0566"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("E35").Value = "  +8.93%  "
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("D37").Value = "'1.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.29%  "
$ws.Range("E38").Value = "  +3.51%  "
$ws.Range("E39").Value = "  +4.50%  "
$ws.Range("D40").Value = "'90.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "'1.356.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.67%  "
$ws.Range("D42").Value = "'15.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.65%  "
$ws.Range("D43").Value = "'13.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +48.33%  "
$ws.Range("D44").Value = "'0.0594"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.81%  "
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  +5.84%  "
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "'45.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +33.62%  "
$ws.Range("D50").Value = "'2.070.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").Value = "'0.0692"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.32%  "
